# Update the cached "datetimeFigureOut" field text (5/7/21 -> 5/10/21) across
# the slide master, every slide layout, and the notes master; then update the
# Q1-Q6 survey-response numbers on the "Survey Responses" slide.
#
# Note: setting TextRange.Text on a range that partially overlaps its old
# value makes the engine emit a minimal-diff set of runs (splitting the run
# at the changed characters). Routing the assignment through an unrelated
# placeholder value first avoids any shared prefix/suffix with the final
# text, so the result collapses back down to a single run, matching how the
# original authors' content was structured.

function Set-CleanText($range, $newText) {
    $range.Text = "zzz__temp_placeholder__zzz"
    $range.Text = $newText
}

$ppPlaceholderDate = 16

$p = $ppt.ActivePresentation

# --- 1. Slide master date placeholder ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        Set-CleanText $sh.TextFrame.TextRange "5/10/21"
    }
}

# --- 2. Every slide layout's date placeholder ---
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            Set-CleanText $sh.TextFrame.TextRange "5/10/21"
        }
    }
}

# --- 3. Notes master date placeholder ---
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $sh = $notesMaster.Shapes.Item($i)
    if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        Set-CleanText $sh.TextFrame.TextRange "5/10/21"
    }
}

# --- 4. Survey responses slide (slide 3): update per-question tallies ---
$slide = $p.Slides.Item(3)
$content = $slide.Shapes.Item(2)
$tr = $content.TextFrame.TextRange

Set-CleanText $tr.Paragraphs(1, 1)  "Q1: 4C, 13D, 1E "
Set-CleanText $tr.Paragraphs(3, 1)  "Q2: 4A, 7B, 6C, 1D "
Set-CleanText $tr.Paragraphs(5, 1)  "Q3: 5A, 10B, 2C, 1D "
Set-CleanText $tr.Paragraphs(7, 1)  "Q4: 10A, 8B "
Set-CleanText $tr.Paragraphs(9, 1)  "Q5: 11A, 4B, 1C, 1D, 1E "
Set-CleanText $tr.Paragraphs(11, 1) "Q6: 5A, 4B, 5C, 2D, 2E "
